$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as other headers (copy style from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:50:34.598038",
    "2021-10-05 10:50:34.598050",
    "2021-10-05 10:50:34.598053",
    "2021-10-05 10:50:34.598056",
    "2021-10-05 10:50:34.598058",
    "2021-10-05 10:50:34.598061",
    "2021-10-05 10:50:34.598064",
    "2021-10-05 10:50:34.598066",
    "2021-10-05 10:50:34.598069",
    "2021-10-05 10:50:34.598072",
    "2021-10-05 10:50:34.598074",
    "2021-10-05 10:50:34.598076",
    "2021-10-05 10:50:34.598079",
    "2021-10-05 10:50:34.598081",
    "2021-10-05 10:50:34.598084",
    "2021-10-05 10:50:34.598086",
    "2021-10-05 10:50:34.598089",
    "2021-10-05 10:50:34.598092",
    "2021-10-05 10:50:34.598094",
    "2021-10-05 10:50:34.598097",
    "2021-10-05 10:50:34.598099",
    "2021-10-05 10:50:34.598102",
    "2021-10-05 10:50:34.598104",
    "2021-10-05 10:50:34.598107",
    "2021-10-05 10:50:34.598109",
    "2021-10-05 10:50:34.598112",
    "2021-10-05 10:50:34.598115",
    "2021-10-05 10:50:34.598117",
    "2021-10-05 10:50:34.598120",
    "2021-10-05 10:50:34.598122",
    "2021-10-05 10:50:34.598125",
    "2021-10-05 10:50:34.598127",
    "2021-10-05 10:50:34.598130",
    "2021-10-05 10:50:34.598132",
    "2021-10-05 10:50:34.598135",
    "2021-10-05 10:50:34.598137",
    "2021-10-05 10:50:34.598140",
    "2021-10-05 10:50:34.598142",
    "2021-10-05 10:50:34.598145",
    "2021-10-05 10:50:34.598148",
    "2021-10-05 10:50:34.598151"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
